$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row that held id=2 (score 52.44); this shifts the old row4
# (id=3, score 82) up to become row 3.
$ws.Rows.Item(3).Delete()

# --- Row 2 (id = 1) gets new score / reasons ---
$ws.Range("B2").Value2 = 83.91
$ws.Range("C2").Value2 = "The candidate has a high score of 83.91 as they have demonstrated strong skills in NLP, Computer Vision, and Python, which are required for the job. Their projects, such as LLMGuard and CloudPhysician's Vital Extraction Challenge, directly relate to the job description, showcasing their expertise in analyzing text data and image processing. The candidate's experience with ReactJS and Django also shows their ability to work with relevant frameworks. Overall, their skills and project experiences align well with the job requirements, justifying the high score."
$ws.Range("D2").Value2 = 4
$ws.Range("E2").Value2 = "The applicant's willingness to come to Japan is evident in their strong interest in AI and the work culture of Japan. They also express adaptability to other cultures and a willingness to learn a new language. Their career plan aligns with the company's focus on AI/ML technologies. The applicant demonstrates good teamwork and problem-solving skills, as well as a willingness to address areas for improvement. Their unique background in both engineering and graphic design, along with specialization in data analytics, sets them apart and adds value to the company. Overall, the applicant's personality and willingness to come to Japan align well with the company's requirements, with slight room for further self-improvement."

# --- Row 3 (id = 3, formerly row 4) gets new score / reasons ---
$ws.Range("A3").Value2 = 3
$ws.Range("B3").Value2 = 83.59
$ws.Range("C3").Value2 = "The candidate's extensive experience in projects related to NLP, computer vision, and Python, along with their skills in NLP, PyTorch, Computer Vision, and Python, make them well-suited for the NLP Engineer Remote position. Their project experience demonstrates expertise in areas required for the job, resulting in a high score of 83.59."
$ws.Range("D3").Value2 = 4
$ws.Range("E3").Value2 = "The applicant demonstrates strong willingness to come to Japan for work, citing interest in the country's advanced AI field and work culture. They have a clear career plan and are open to working in different countries, showcasing adaptability. Their expectation from the company aligns with their career goals, and they express a desire to contribute to cutting-edge technologies. They also show an understanding of their strengths, weaknesses, and have taken steps to address their areas for improvement. Their experience in problem-solving and teamwork, along with a background in both engineering and graphic design, positions them uniquely for the role. However, they acknowledge the need to improve their solo work skills, which is an important aspect for the company's requirement of adaptability. Hence, a score of 4 is given. "
